# "Wireframes version 2." -> "Wireframes version 1."
# (revert of a revert: the heading text goes from "Version 2." back to "Version 1.")
$d = $word.ActiveDocument

# 1. Turn the " 2" run into " 1." (stays inside its own run, doesn't touch the
#    bookmark that sits right after it).
$d.Content.Find.Execute("2", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1.", 2)

# 2. The text is now "Version 1.." - the trailing "." is the old standalone
#    run that used to close the sentence. Remove just that last character so
#    the bookmark (which sits before it) is left untouched.
#    Note: Content.Text.Length includes the implicit end-of-story mark, so
#    the last real character sits at [len-2, len-1).
$len = $d.Content.Text.Length
$trailing = $d.Range($len - 2, $len - 1)
$trailing.Delete()

# 3. Force Word to normalize/merge the "Versi" + "on" runs (identical,
#    formatting-less runs) into a single "Version" run, matching how Word
#    coalesces runs after an in-place replace.
$d.Content.Find.Execute("Versi", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Versi", 2)
